# ActualizarDatosSeguridad.xlsx - "Actualizacion data y funcionalidades"
#
# The commit only changes two data values on the "Datos" sheet (row 2):
#   - D2 (usuario):  autotest32               -> autotest10
#   - M2 (correo):   jfernandez@todo1.net     -> automatizaciontodo1@gmail.com
#     (M2 keeps its existing hyperlink formatting/link; only the displayed
#      text changes, exactly like typing a new value over the old one.)
#
# Everything else in the published diff (fileVersion/xr revision bookkeeping,
# column-width micro-adjustments, cellXfs/font re-numbering after Excel
# recomputes its style table, selection/scroll position, etc.) is
# incidental save noise produced by the newer Excel build that resaved the
# workbook - it carries no data/formatting intent of its own, so it is not
# reproduced here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Write M2 first, then D2, so the shared-string table grows in the same
# order the workbook was edited in (email added, then the usuario tweak).
$ws.Range("M2").Value = "automatizaciontodo1@gmail.com"
$ws.Range("D2").Value = "autotest10"

# Restore the cursor position recorded in the saved file.
$ws.Range("E9").Select()
